# Updated cryptos list on Thu Apr 20 06:43:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.108.92"
$ws.Range("E2").Value = "  -3.83%  "

$ws.Range("D3").Value = "1.967.56"
$ws.Range("E3").Value = "  -5.75%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.81%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("E7").Value = "  -5.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4219"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.09%  "

$ws.Range("D13").Value = "2.015.84"
$ws.Range("E13").Value = "  -5.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.880"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.446"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -9.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06687"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.979"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.99%  "

$ws.Range("D23").Value = "29.139.42"
$ws.Range("E23").Value = "  -3.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.287"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("D26").Value = "2.250.29"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.198"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.265"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.046"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09862"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.530"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.790"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.681"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02432"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.049"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.302"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06368"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6462"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1990"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.72%  "

$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6247"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.198"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.289"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.473"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000333"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06981"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "

